# Update LR-pair TPM-derived metrics (Ptn-Plxnb2) to reflect recalculated values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.2563003333333334
$ws.Range("H2").Value = 0.7689010000000001
$ws.Range("I2").Value = 0.02986826554325775
$ws.Range("J2").Value = 0.02986826554325775
$ws.Range("M2").Value = 12.046506
$ws.Range("N2").Value = 36.139518
$ws.Range("O2").Value = 0.06409289449370618
$ws.Range("P2").Value = 0.06409289449370618
$ws.Range("Q2").Value = 3.087523503302001
$ws.Range("R2").Value = 27.787711529718
$ws.Range("S2").Value = 0.001914343592174019
$ws.Range("T2").Value = 0.001914343592174019

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.2563003333333334
$ws.Range("H3").Value = 0.7689010000000001
$ws.Range("I3").Value = 0.02986826554325775
$ws.Range("J3").Value = 0.02986826554325775
$ws.Range("O3").Value = 0.199479586067244
$ws.Range("P3").Value = 0.199479586067244
$ws.Range("Q3").Value = 9.609456949585113
$ws.Range("R3").Value = 86.48511254626601
$ws.Range("S3").Value = 0.005958109247115584
$ws.Range("T3").Value = 0.005958109247115583

$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2563003333333334
$ws.Range("H4").Value = 0.7689010000000001
$ws.Range("I4").Value = 0.02986826554325775
$ws.Range("J4").Value = 0.02986826554325775
$ws.Range("M4").Value = 53.23753466666667
$ws.Range("N4").Value = 159.712604
$ws.Range("O4").Value = 0.2832479137515634
$ws.Range("P4").Value = 0.2832479137515634
$ws.Range("Q4").Value = 13.64479788091156
$ws.Range("R4").Value = 122.803180928204
$ws.Range("S4").Value = 0.008460123902505464
$ws.Range("T4").Value = 0.008460123902505464

$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2563003333333334
$ws.Range("H5").Value = 0.7689010000000001
$ws.Range("I5").Value = 0.02986826554325775
$ws.Range("J5").Value = 0.02986826554325775
$ws.Range("M5").Value = 15.283152
$ws.Range("N5").Value = 45.849456
$ws.Range("O5").Value = 0.08131332426740788
$ws.Range("P5").Value = 0.08131332426740788
$ws.Range("Q5").Value = 3.917076951984001
$ws.Range("R5").Value = 35.25369256785601
$ws.Range("S5").Value = 0.002428687961423963
$ws.Range("T5").Value = 0.002428687961423963

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2563003333333334
$ws.Range("H6").Value = 0.7689010000000001
$ws.Range("I6").Value = 0.02986826554325775
$ws.Range("J6").Value = 0.02986826554325775
$ws.Range("M6").Value = 10.028695
$ws.Range("N6").Value = 30.086085
$ws.Range("O6").Value = 0.05335722163294142
$ws.Range("P6").Value = 0.05335722163294142
$ws.Range("Q6").Value = 2.570357871398334
$ws.Range("R6").Value = 23.133220842585
$ws.Range("S6").Value = 0.001593687664383151
$ws.Range("T6").Value = 0.001593687664383151

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2563003333333334
$ws.Range("H7").Value = 0.7689010000000001
$ws.Range("I7").Value = 0.02986826554325775
$ws.Range("J7").Value = 0.02986826554325775
$ws.Range("M7").Value = 59.86500266666667
$ws.Range("N7").Value = 179.595008
$ws.Range("O7").Value = 0.3185090597871371
$ws.Range("P7").Value = 0.3185090597871371
$ws.Range("Q7").Value = 15.34342013846756
$ws.Range("R7").Value = 138.090781246208
$ws.Range("S7").Value = 0.009513313175655569
$ws.Range("T7").Value = 0.009513313175655567

$ws.Range("G8").Value = 6.495645000000001
$ws.Range("I8").Value = 0.7569777503270297
$ws.Range("J8").Value = 0.7569777503270296
$ws.Range("M8").Value = 12.046506
$ws.Range("N8").Value = 36.139518
$ws.Range("O8").Value = 0.06409289449370618
$ws.Range("P8").Value = 0.06409289449370618
$ws.Range("Q8").Value = 78.24982646637001
$ws.Range("R8").Value = 704.2484381973302
$ws.Range("S8").Value = 0.04851689508579338
$ws.Range("T8").Value = 0.04851689508579337

$ws.Range("G9").Value = 6.495645000000001
$ws.Range("I9").Value = 0.7569777503270297
$ws.Range("J9").Value = 0.7569777503270296
$ws.Range("O9").Value = 0.199479586067244
$ws.Range("P9").Value = 0.199479586067244
$ws.Range("S9").Value = 0.1510016082973495
$ws.Range("T9").Value = 0.1510016082973495

$ws.Range("G10").Value = 6.495645000000001
$ws.Range("I10").Value = 0.7569777503270297
$ws.Range("J10").Value = 0.7569777503270296
$ws.Range("M10").Value = 53.23753466666667
$ws.Range("N10").Value = 159.712604
$ws.Range("O10").Value = 0.2832479137515634
$ws.Range("P10").Value = 0.2832479137515634
$ws.Range("Q10").Value = 345.81212586986
$ws.Range("R10").Value = 3112.309132828741
$ws.Range("S10").Value = 0.214412368536483
$ws.Range("T10").Value = 0.214412368536483

$ws.Range("G11").Value = 6.495645000000001
$ws.Range("I11").Value = 0.7569777503270297
$ws.Range("J11").Value = 0.7569777503270296
$ws.Range("M11").Value = 15.283152
$ws.Range("N11").Value = 45.849456
$ws.Range("O11").Value = 0.08131332426740788
$ws.Range("P11").Value = 0.08131332426740788
$ws.Range("Q11").Value = 99.27392987304002
$ws.Range("R11").Value = 893.4653688573602
$ws.Range("S11").Value = 0.06155237727555469
$ws.Range("T11").Value = 0.06155237727555468

$ws.Range("G12").Value = 6.495645000000001
$ws.Range("I12").Value = 0.7569777503270297
$ws.Range("J12").Value = 0.7569777503270296
$ws.Range("M12").Value = 10.028695
$ws.Range("N12").Value = 30.086085
$ws.Range("O12").Value = 0.05335722163294142
$ws.Range("P12").Value = 0.05335722163294142
$ws.Range("Q12").Value = 65.14284253327502
$ws.Range("R12").Value = 586.285582799475
$ws.Range("S12").Value = 0.04039022959540473
$ws.Range("T12").Value = 0.04039022959540472

$ws.Range("G13").Value = 6.495645000000001
$ws.Range("I13").Value = 0.7569777503270297
$ws.Range("J13").Value = 0.7569777503270296
$ws.Range("M13").Value = 59.86500266666667
$ws.Range("N13").Value = 179.595008
$ws.Range("O13").Value = 0.3185090597871371
$ws.Range("P13").Value = 0.3185090597871371
$ws.Range("Q13").Value = 388.8618052467201
$ws.Range("R13").Value = 3499.75624722048
$ws.Range("S13").Value = 0.2411042715364444
$ws.Range("T13").Value = 0.2411042715364444

$ws.Range("G14").Value = 1.804372666666667
$ws.Range("H14").Value = 5.413118000000001
$ws.Range("I14").Value = 0.2102747243676212
$ws.Range("J14").Value = 0.2102747243676212
$ws.Range("M14").Value = 12.046506
$ws.Range("N14").Value = 36.139518
$ws.Range("O14").Value = 0.06409289449370618
$ws.Range("P14").Value = 0.06409289449370618
$ws.Range("Q14").Value = 21.736386155236
$ws.Range("R14").Value = 195.627475397124
$ws.Range("S14").Value = 0.01347711572358709
$ws.Range("T14").Value = 0.01347711572358709

$ws.Range("G15").Value = 1.804372666666667
$ws.Range("H15").Value = 5.413118000000001
$ws.Range("I15").Value = 0.2102747243676212
$ws.Range("J15").Value = 0.2102747243676212
$ws.Range("O15").Value = 0.199479586067244
$ws.Range("P15").Value = 0.199479586067244
$ws.Range("Q15").Value = 67.6512637960209
$ws.Range("R15").Value = 608.8613741641881
$ws.Range("S15").Value = 0.04194551497725691
$ws.Range("T15").Value = 0.0419455149772569

$ws.Range("G16").Value = 1.804372666666667
$ws.Range("H16").Value = 5.413118000000001
$ws.Range("I16").Value = 0.2102747243676212
$ws.Range("J16").Value = 0.2102747243676212
$ws.Range("M16").Value = 53.23753466666667
$ws.Range("N16").Value = 159.712604
$ws.Range("O16").Value = 0.2832479137515634
$ws.Range("P16").Value = 0.2832479137515634
$ws.Range("Q16").Value = 96.06035239325246
$ws.Range("R16").Value = 864.5431715392721
$ws.Range("S16").Value = 0.05955987699181375
$ws.Range("T16").Value = 0.05955987699181374

$ws.Range("G17").Value = 1.804372666666667
$ws.Range("H17").Value = 5.413118000000001
$ws.Range("I17").Value = 0.2102747243676212
$ws.Range("J17").Value = 0.2102747243676212
$ws.Range("M17").Value = 15.283152
$ws.Range("N17").Value = 45.849456
$ws.Range("O17").Value = 0.08131332426740788
$ws.Range("P17").Value = 0.08131332426740788
$ws.Range("Q17").Value = 27.576501729312
$ws.Range("R17").Value = 248.1885155638081
$ws.Range("S17").Value = 0.0170981368477442
$ws.Range("T17").Value = 0.01709813684774419

$ws.Range("G18").Value = 1.804372666666667
$ws.Range("H18").Value = 5.413118000000001
$ws.Range("I18").Value = 0.2102747243676212
$ws.Range("J18").Value = 0.2102747243676212
$ws.Range("M18").Value = 10.028695
$ws.Range("N18").Value = 30.086085
$ws.Range("O18").Value = 0.05335722163294142
$ws.Range("P18").Value = 0.05335722163294142
$ws.Range("Q18").Value = 18.09550314033667
$ws.Range("R18").Value = 162.85952826303
$ws.Range("S18").Value = 0.01121967507188883
$ws.Range("T18").Value = 0.01121967507188883

$ws.Range("G19").Value = 1.804372666666667
$ws.Range("H19").Value = 5.413118000000001
$ws.Range("I19").Value = 0.2102747243676212
$ws.Range("J19").Value = 0.2102747243676212
$ws.Range("M19").Value = 59.86500266666667
$ws.Range("N19").Value = 179.595008
$ws.Range("O19").Value = 0.3185090597871371
$ws.Range("P19").Value = 0.3185090597871371
$ws.Range("Q19").Value = 108.0187745016605
$ws.Range("R19").Value = 972.1689705149441
$ws.Range("S19").Value = 0.06697440475533042
$ws.Range("T19").Value = 0.06697440475533042

$ws.Range("I20").Value = 0.002879259762091359
$ws.Range("J20").Value = 0.002879259762091358
$ws.Range("M20").Value = 12.046506
$ws.Range("N20").Value = 36.139518
$ws.Range("O20").Value = 0.06409289449370618
$ws.Range("P20").Value = 0.06409289449370618
$ws.Range("Q20").Value = 0.2976330237420001
$ws.Range("R20").Value = 2.678697213678
$ws.Range("S20").Value = 0.000184540092151695
$ws.Range("T20").Value = 0.000184540092151695

$ws.Range("I21").Value = 0.002879259762091359
$ws.Range("J21").Value = 0.002879259762091358
$ws.Range("O21").Value = 0.199479586067244
$ws.Range("P21").Value = 0.199479586067244
$ws.Range("S21").Value = 0.0005743535455220558
$ws.Range("T21").Value = 0.0005743535455220557

$ws.Range("I22").Value = 0.002879259762091359
$ws.Range("J22").Value = 0.002879259762091358
$ws.Range("M22").Value = 53.23753466666667
$ws.Range("N22").Value = 159.712604
$ws.Range("O22").Value = 0.2832479137515634
$ws.Range("P22").Value = 0.2832479137515634
$ws.Range("Q22").Value = 1.315339769009334
$ws.Range("R22").Value = 11.838057921084
$ws.Range("S22").Value = 0.0008155443207612002
$ws.Range("T22").Value = 0.0008155443207612001

$ws.Range("I23").Value = 0.002879259762091359
$ws.Range("J23").Value = 0.002879259762091358
$ws.Range("M23").Value = 15.283152
$ws.Range("N23").Value = 45.849456
$ws.Range("O23").Value = 0.08131332426740788
$ws.Range("P23").Value = 0.08131332426740788
$ws.Range("Q23").Value = 0.3776008364640001
$ws.Range("R23").Value = 3.398407528176
$ws.Range("S23").Value = 0.0002341221826850343
$ws.Range("T23").Value = 0.0002341221826850343

$ws.Range("I24").Value = 0.002879259762091359
$ws.Range("J24").Value = 0.002879259762091358
$ws.Range("M24").Value = 10.028695
$ws.Range("N24").Value = 30.086085
$ws.Range("O24").Value = 0.05335722163294142
$ws.Range("P24").Value = 0.05335722163294142
$ws.Range("Q24").Value = 0.247778967365
$ws.Range("R24").Value = 2.230010706285
$ws.Range("S24").Value = 0.0001536293012647188
$ws.Range("T24").Value = 0.0001536293012647188

$ws.Range("I25").Value = 0.002879259762091359
$ws.Range("J25").Value = 0.002879259762091358
$ws.Range("M25").Value = 59.86500266666667
$ws.Range("N25").Value = 179.595008
$ws.Range("O25").Value = 0.3185090597871371
$ws.Range("P25").Value = 0.3185090597871371
$ws.Range("Q25").Value = 1.479084620885334
$ws.Range("R25").Value = 13.311761587968
$ws.Range("S25").Value = 0.0009170703197066546
$ws.Range("T25").Value = 0.0009170703197066545
